$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '69.704.08'
$ws.Cells.Item(2, 5).Value = '  +2.64%  '
$ws.Cells.Item(3, 4).Value = '3.411.53'
$ws.Cells.Item(3, 5).Value = '  +2.14%  '
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.NumberFormat = "General"
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '587.58'
$c.NumberFormat = "General"
$ws.Cells.Item(5, 5).Value = '  +0.75%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '181.94'
$c.NumberFormat = "General"
$ws.Cells.Item(6, 5).Value = '  +3.91%  '
$ws.Cells.Item(7, 5).Value = '  +1.89%  '
$ws.Cells.Item(8, 5).Value = '  +0.03%  '
$ws.Cells.Item(9, 5).Value = '  +11.50%  '
$ws.Cells.Item(10, 5).Value = '  +2.52%  '
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '48.60'
$c.NumberFormat = "General"
$ws.Cells.Item(11, 5).Value = '  +2.81%  '
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '0.0000289'
$c.NumberFormat = "General"
$ws.Cells.Item(12, 5).Value = '  +5.20%  '
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '691.02'
$c.NumberFormat = "General"
$ws.Cells.Item(13, 5).Value = '  -0.95%  '
$ws.Cells.Item(14, 5).Value = '  +4.31%  '
$ws.Cells.Item(15, 4).Value = '3.963.45'
$ws.Cells.Item(15, 5).Value = '  +2.22%  '
$ws.Cells.Item(16, 4).Value = '69.715.13'
$ws.Cells.Item(16, 5).Value = '  +2.60%  '
$ws.Cells.Item(17, 5).Value = '  +1.66%  '
$ws.Cells.Item(18, 4).Value = '3.407.07'
$ws.Cells.Item(18, 5).Value = '  +2.17%  '
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '17.82'
$c.NumberFormat = "General"
$ws.Cells.Item(19, 5).Value = '  +2.22%  '
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '11.41'
$c.NumberFormat = "General"
$ws.Cells.Item(20, 5).Value = '  +2.49%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '17.39'
$c.NumberFormat = "General"
$ws.Cells.Item(22, 5).Value = '  +2.42%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '5.38'
$c.NumberFormat = "General"
$ws.Cells.Item(23, 5).Value = '  -0.68%  '
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '103.56'
$c.NumberFormat = "General"
$ws.Cells.Item(24, 5).Value = '  +2.29%  '
$ws.Cells.Item(25, 5).Value = '  +1.26%  '
$ws.Cells.Item(26, 5).Value = '  +2.06%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '9.78'
$c.NumberFormat = "General"
$ws.Cells.Item(27, 5).Value = '  +3.66%  '
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '34.15'
$c.NumberFormat = "General"
$ws.Cells.Item(28, 5).Value = '  +3.29%  '
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '8.88'
$c.NumberFormat = "General"
$ws.Cells.Item(29, 5).Value = '  +4.07%  '
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '6.99'
$c.NumberFormat = "General"
$ws.Cells.Item(30, 5).Value = '  +0.45%  '
$ws.Cells.Item(31, 2).Value = 'dogwifhat'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '3.70'
$c.NumberFormat = "General"
$ws.Cells.Item(31, 5).Value = '  +11.11%  '
$ws.Cells.Item(32, 2).Value = 'Cosmos'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '11.19'
$c.NumberFormat = "General"
$ws.Cells.Item(32, 5).Value = '  +1.61%  '
$ws.Cells.Item(33, 2).Value = 'Bittensor'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '558.44'
$c.NumberFormat = "General"
$ws.Cells.Item(33, 5).Value = '  -2.87%  '
$ws.Cells.Item(34, 5).Value = '  +1.69%  '
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '58.73'
$c.NumberFormat = "General"
$ws.Cells.Item(35, 5).Value = '  +3.60%  '
$ws.Cells.Item(36, 5).Value = '  -0.02%  '
$ws.Cells.Item(37, 4).Value = '3.664.77'
$ws.Cells.Item(37, 5).Value = '  -2.53%  '
$ws.Cells.Item(38, 5).Value = '  +5.73%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '36.04'
$c.NumberFormat = "General"
$ws.Cells.Item(39, 5).Value = '  +1.69%  '
$ws.Cells.Item(40, 4).Value = '0.0₃0741'
$ws.Cells.Item(40, 5).Value = '  +9.39%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '3.29'
$c.NumberFormat = "General"
$ws.Cells.Item(41, 5).Value = '  +4.42%  '
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '2.70'
$c.NumberFormat = "General"
$ws.Cells.Item(42, 5).Value = '  +3.57%  '
$ws.Cells.Item(43, 5).Value = '  +2.21%  '
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.0430'
$c.NumberFormat = "General"
$ws.Cells.Item(44, 5).Value = '  +5.98%  '
$ws.Cells.Item(45, 5).Value = '  +1.09%  '
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '2.69'
$c.NumberFormat = "General"
$ws.Cells.Item(46, 5).Value = '  +1.97%  '
$ws.Cells.Item(47, 5).Value = '  +1.15%  '
$ws.Cells.Item(48, 5).Value = '  +4.77%  '
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.NumberFormat = "General"
$ws.Cells.Item(49, 5).Value = '  -0.30%  '
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '131.28'
$c.NumberFormat = "General"
$ws.Cells.Item(50, 5).Value = '  +0.95%  '
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '2.72'
$c.NumberFormat = "General"
$ws.Cells.Item(51, 5).Value = '  +0.80%  '
